$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 (IMEI 359998100056409 / "B") was removed from the sales sheet.
# Clear its contents and reset the row height back to the sheet default so
# that the now-empty row is dropped entirely (rather than leaving a blank
# placeholder), matching the rows below (13, 15, 17, ...) which keep their
# original row numbers.
$ws.Range("A11:B11").Clear()
$ws.Rows("11:11").AutoFit()

# Restore the workbook window position/size recorded at last save.
$win = $excel.Windows.Item(1)
$win.WindowState = -4143  # xlNormal
$win.Left = 1380
$win.Top = 510
$win.Width = 11910
$win.Height = 10560
